$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1575
$ws1.Range("F4").Value = 33

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1575
$ws4.Range("F4").Value = 33
